$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B2").Value = '112号直流'
$ws.Range("C2").Value = 46037.393877314818
$ws.Range("D2").Value = 46038.289571759262

$ws.Range("A3").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B3").Value = '210号直流'
$ws.Range("C3").Value = 46037.453668981485
$ws.Range("D3").Value = 46038.289571759262

$ws.Range("A4").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B4").Value = '106号直流'
$ws.Range("C4").Value = 46037.552245370367
$ws.Range("D4").Value = 46038.289571759262

$ws.Range("A5").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B5").Value = '111号直流'
$ws.Range("C5").Value = 46037.57135416667
$ws.Range("D5").Value = 46038.289571759262

$ws.Range("A6").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B6").Value = '109号直流'
$ws.Range("C6").Value = 46037.623981481483
$ws.Range("D6").Value = 46038.289571759262

$ws.Range("A7").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B7").Value = '204号直流'
$ws.Range("C7").Value = 46037.664525462962
$ws.Range("D7").Value = 46038.289571759262

$ws.Range("A8").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B8").Value = '305号直流'
$ws.Range("C8").Value = 46037.706655092596
$ws.Range("D8").Value = 46038.289571759262

$ws.Range("A9").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B9").Value = '206号直流'
$ws.Range("C9").Value = 46037.72729166667
$ws.Range("D9").Value = 46038.289571759262

$ws.Range("A10").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B10").Value = '309号直流'
$ws.Range("C10").Value = 46037.740208333336
$ws.Range("D10").Value = 46038.289571759262

$ws.Range("A11").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B11").Value = '108号直流'
$ws.Range("C11").Value = 46037.745648148149
$ws.Range("D11").Value = 46038.289571759262

$ws.Range("A12").Value = '飞狐四方坪西区充电站'
$ws.Range("B12").Value = '9176699400500604'
$ws.Range("C12").Value = 46036.551874999997
$ws.Range("D12").Value = 46038.303101851852

$ws.Range("A13").Value = '飞狐四方坪西区充电站'
$ws.Range("B13").Value = '9176699400501205'
$ws.Range("C13").Value = 46036.768391203703
$ws.Range("D13").Value = 46038.303101851852

$ws.Range("A14").Value = '飞狐四方坪西区充电站'
$ws.Range("B14").Value = '9176699400501202'
$ws.Range("C14").Value = 46037.046331018515
$ws.Range("D14").Value = 46038.303101851852

$ws.Range("A15").Value = '飞狐四方坪西区充电站'
$ws.Range("B15").Value = '9176699400500904'
$ws.Range("C15").Value = 46037.079467592594
$ws.Range("D15").Value = 46038.303101851852

$ws.Range("A16").Value = '飞狐四方坪西区充电站'
$ws.Range("B16").Value = '9176699400500203'
$ws.Range("C16").Value = 46037.146307870367
$ws.Range("D16").Value = 46038.303101851852

$ws.Range("A17").Value = '飞狐四方坪西区充电站'
$ws.Range("B17").Value = '9176699400501201'
$ws.Range("C17").Value = 46037.193379629629
$ws.Range("D17").Value = 46038.303101851852

$ws.Range("A18").Value = '飞狐四方坪西区充电站'
$ws.Range("B18").Value = '9176699400500304'
$ws.Range("C18").Value = 46037.196342592593
$ws.Range("D18").Value = 46038.303101851852

$ws.Range("A19").Value = '飞狐四方坪西区充电站'
$ws.Range("B19").Value = '9176699400500201'
$ws.Range("C19").Value = 46037.19939814815
$ws.Range("D19").Value = 46038.303101851852

$ws.Range("A20").Value = '飞狐四方坪南区充电站'
$ws.Range("B20").Value = '9176699368200406'
$ws.Range("C20").Value = 46037.328738425924
$ws.Range("D20").Value = 46038.303101851852

$ws.Range("A21").Value = '飞狐四方坪西区充电站'
$ws.Range("B21").Value = '9176699400501302'
$ws.Range("C21").Value = 46037.526284722226
$ws.Range("D21").Value = 46038.303101851852

$ws.Range("A22").Value = '飞狐四方坪东区充电站'
$ws.Range("B22").Value = '9176699416300203'
$ws.Range("C22").Value = 46037.531782407408
$ws.Range("D22").Value = 46038.303101851852

$ws.Range("A23").Value = '飞狐四方坪西区充电站'
$ws.Range("B23").Value = '9176699400500105'
$ws.Range("C23").Value = 46037.537627314814
$ws.Range("D23").Value = 46038.303101851852

$ws.Range("A24").Value = '飞狐四方坪东区充电站'
$ws.Range("B24").Value = '9176699442100402'
$ws.Range("C24").Value = 46037.546168981484
$ws.Range("D24").Value = 46038.303101851852

$ws.Range("A25").Value = '飞狐四方坪南区充电站'
$ws.Range("B25").Value = '9176699368200103'
$ws.Range("C25").Value = 46037.548171296294
$ws.Range("D25").Value = 46038.303101851852

$ws.Range("A26").Value = '飞狐四方坪西区充电站'
$ws.Range("B26").Value = '9176699400501303'
$ws.Range("C26").Value = 46037.548541666663
$ws.Range("D26").Value = 46038.303101851852

$ws.Range("A27").Value = '飞狐四方坪西区充电站'
$ws.Range("B27").Value = '9176699400500303'
$ws.Range("C27").Value = 46037.552986111114
$ws.Range("D27").Value = 46038.303101851852

$ws.Range("A28").Value = '飞狐四方坪西区充电站'
$ws.Range("B28").Value = '9176699400501203'
$ws.Range("C28").Value = 46037.553449074076
$ws.Range("D28").Value = 46038.303101851852

$ws.Range("A29").Value = '飞狐四方坪南区充电站'
$ws.Range("B29").Value = '9176699368200203'
$ws.Range("C29").Value = 46037.55810185185
$ws.Range("D29").Value = 46038.303101851852

$ws.Range("A30").Value = '飞狐四方坪东区充电站'
$ws.Range("B30").Value = '9176699442100302'
$ws.Range("C30").Value = 46037.567824074074
$ws.Range("D30").Value = 46038.303101851852

$ws.Range("A31").Value = '飞狐四方坪西区充电站'
$ws.Range("B31").Value = '9176699400500501'
$ws.Range("C31").Value = 46037.569108796299
$ws.Range("D31").Value = 46038.303101851852

$ws.Range("A32").Value = '飞狐四方坪西区充电站'
$ws.Range("B32").Value = '9176699400500605'
$ws.Range("C32").Value = 46037.571388888886
$ws.Range("D32").Value = 46038.303101851852

$ws.Range("A33").Value = '飞狐四方坪东区充电站'
$ws.Range("B33").Value = '9176699425700302'
$ws.Range("C33").Value = 46037.572291666664
$ws.Range("D33").Value = 46038.303101851852

$ws.Range("A34").Value = '飞狐四方坪西区充电站'
$ws.Range("B34").Value = '9176699400501204'
$ws.Range("C34").Value = 46037.6096875
$ws.Range("D34").Value = 46038.303101851852

$ws.Range("A35").Value = '飞狐四方坪西区充电站'
$ws.Range("B35").Value = '9176699400500102'
$ws.Range("C35").Value = 46037.61037037037
$ws.Range("D35").Value = 46038.303101851852

$ws.Range("A36").Value = '飞狐四方坪西区充电站'
$ws.Range("B36").Value = '9176699400501104'
$ws.Range("C36").Value = 46037.620081018518
$ws.Range("D36").Value = 46038.303101851852

$ws.Range("A37").Value = '飞狐四方坪西区充电站'
$ws.Range("B37").Value = '9176699400501304'
$ws.Range("C37").Value = 46037.635069444441
$ws.Range("D37").Value = 46038.303101851852

$ws.Range("A38").Value = '飞狐四方坪南区充电站'
$ws.Range("B38").Value = '9176699368200305'
$ws.Range("C38").Value = 46037.694247685184
$ws.Range("D38").Value = 46038.303101851852

$ws.Range("A39").Value = '飞狐四方坪西区充电站'
$ws.Range("B39").Value = '9176699400501101'
$ws.Range("C39").Value = 46037.695520833331
$ws.Range("D39").Value = 46038.303101851852

$ws.Range("A40").Value = '飞狐四方坪西区充电站'
$ws.Range("B40").Value = '9176699400501105'
$ws.Range("C40").Value = 46037.701458333337
$ws.Range("D40").Value = 46038.303101851852

$ws.Range("A41").Value = '飞狐四方坪南区充电站'
$ws.Range("B41").Value = '9176699368200304'
$ws.Range("C41").Value = 46037.704907407409
$ws.Range("D41").Value = 46038.303101851852

$ws.Range("A42").Value = '飞狐四方坪南区充电站'
$ws.Range("B42").Value = '9176699368200201'
$ws.Range("C42").Value = 46037.706458333334
$ws.Range("D42").Value = 46038.303101851852

$ws.Range("A43").Value = '飞狐四方坪西区充电站'
$ws.Range("B43").Value = '9176699400500302'
$ws.Range("C43").Value = 46037.755335648151
$ws.Range("D43").Value = 46038.303101851852

$ws.Range("A44").Value = '飞狐四方坪西区充电站'
$ws.Range("B44").Value = '9176699400500205'
$ws.Range("C44").Value = 46037.772245370368
$ws.Range("D44").Value = 46038.303101851852

$ws.Range("E17").Select()
